$d = $word.ActiveDocument

# --- Name of the invitee: Арсению Эдуардовичу Губанову -> Андрею Александровичу Зелиховскому ---
$d.Content.Find.Execute("Арсению", $true, $false, $false, $false, $false, $true, 1, $false, "Андрею", 2) | Out-Null
$d.Content.Find.Execute("Эдуардовичу", $true, $false, $false, $false, $false, $true, 1, $false, "Александровичу", 2) | Out-Null
$d.Content.Find.Execute("Губанову", $true, $false, $false, $false, $false, $true, 1, $false, "Зелиховскому", 2) | Out-Null

# --- Trailing underscores right after the name, before the {cDate} field shrink from 8 to 5 ---
$d.Content.Find.Execute("в том, что________", $true, $false, $false, $false, $false, $true, 1, $false, "в том, что_____", 2) | Out-Null

# --- Underscore runs after {cDate} field ---
$d.Content.Find.Execute("187264", $true, $false, $false, $false, $false, $true, 1, $false, "756628", 2) | Out-Null
